$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, shifting existing rows 46..114 down to 47..115.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with its data.
$ws.Cells.Item(46, 1).Value = 8
$ws.Cells.Item(46, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44664
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100109
$ws.Cells.Item(46, 8).Value = "Uva"
$ws.Cells.Item(46, 9).Value = 100109001
$ws.Cells.Item(46, 10).Value = "Uva"
$ws.Cells.Item(46, 11).Value = "Red Globe"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 400
$ws.Cells.Item(46, 14).Value = 9500
$ws.Cells.Item(46, 15).Value = 10000
$ws.Cells.Item(46, 16).Value = 9750
$ws.Cells.Item(46, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(46, 19).Value = 542
$ws.Cells.Item(46, 20).Value = 18
